$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices & 1h volume %) plus two coin-row
# reorderings, matching the upstream scrape refresh.

$ws.Range("D2").Value2 = "64.349.20"
$ws.Range("E2").Value2 = "  +0.01%  "
$ws.Range("D3").Value2 = "3.481.59"
$ws.Range("E3").Value2 = "  -0.15%  "
$ws.Range("E4").Value2 = "  -0.02%  "
$ws.Range("D5").Value2 = "'585.39"
$ws.Range("E5").Value2 = "  +0.42%  "
$ws.Range("D6").Value2 = "'134.05"
$ws.Range("E6").Value2 = "  +1.87%  "
$ws.Range("D7").Value2 = "3.482.87"
$ws.Range("E7").Value2 = "  -0.15%  "
$ws.Range("E8").Value2 = "  -0.01%  "
$ws.Range("D9").Value2 = "'0.485"
$ws.Range("E9").Value2 = "  -1.08%  "
$ws.Range("E10").Value2 = "  -0.25%  "
$ws.Range("D11").Value2 = "'7.17"
$ws.Range("E11").Value2 = "  -0.90%  "
$ws.Range("E12").Value2 = "  -2.91%  "
$ws.Range("D13").Value2 = "4.075.57"
$ws.Range("E13").Value2 = "  -0.07%  "
$ws.Range("E14").Value2 = "  +1.62%  "
$ws.Range("D15").Value2 = "'0.0000179"
$ws.Range("E15").Value2 = "  +0.54%  "
$ws.Range("D16").Value2 = "3.484.37"
$ws.Range("E16").Value2 = "  +0.00%  "
$ws.Range("D17").Value2 = "64.345.71"
$ws.Range("E17").Value2 = "  -0.05%  "
$ws.Range("D18").Value2 = "'25.05"
$ws.Range("E18").Value2 = "  -9.76%  "
$ws.Range("D19").Value2 = "'9.94"
$ws.Range("E19").Value2 = "  -0.04%  "
$ws.Range("E20").Value2 = "  -0.15%  "
$ws.Range("D21").Value2 = "'13.69"
$ws.Range("E21").Value2 = "  -4.19%  "
$ws.Range("D22").Value2 = "'385.50"
$ws.Range("E22").Value2 = "  -1.90%  "
$ws.Range("D23").Value2 = "'0.564"
$ws.Range("E23").Value2 = "  -2.02%  "
$ws.Range("D24").Value2 = "3.621.69"
$ws.Range("E24").Value2 = "  -0.11%  "
$ws.Range("E25").Value2 = "  +1.16%  "
$ws.Range("E26").Value2 = "  +0.06%  "
$ws.Range("E27").Value2 = "  +3.94%  "
$ws.Range("E28").Value2 = "  +0.09%  "
$ws.Range("B29").Value2 = "Fetch.AI"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value2 = "'1.53"
$ws.Range("E29").Value2 = "  -1.13%  "
$ws.Range("B30").Value2 = "RenderToken"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value2 = "'7.38"
$ws.Range("E30").Value2 = "  -0.78%  "
$ws.Range("D31").Value2 = "'2.22"
$ws.Range("E31").Value2 = "  -1.09%  "
$ws.Range("D32").Value2 = "'8.17"
$ws.Range("E32").Value2 = "  +0.03%  "
$ws.Range("D33").Value2 = "3.500.57"
$ws.Range("E33").Value2 = "  +0.42%  "
$ws.Range("E34").Value2 = "  -0.03%  "
$ws.Range("E35").Value2 = "  +1.39%  "
$ws.Range("D36").Value2 = "'23.33"
$ws.Range("E36").Value2 = "  -2.29%  "
$ws.Range("D37").Value2 = "'5.29"
$ws.Range("E37").Value2 = "  +1.07%  "
$ws.Range("E38").Value2 = "  -2.12%  "
$ws.Range("D39").Value2 = "'1.53"
$ws.Range("E39").Value2 = "  -2.34%  "
$ws.Range("D40").Value2 = "'162.55"
$ws.Range("E40").Value2 = "  -4.45%  "
$ws.Range("D41").Value2 = "'0.0776"
$ws.Range("E41").Value2 = "  -3.23%  "
$ws.Range("D42").Value2 = "'0.801"
$ws.Range("E42").Value2 = "  -1.27%  "
$ws.Range("B43").Value2 = "FirstDigitalUSD"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value2 = "'1.00"
$ws.Range("E43").Value2 = "  +0.01%  "
$ws.Range("B44").Value2 = "EnergySwap"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value2 = "'25.41"
$ws.Range("E44").Value2 = "  -0.43%  "
$ws.Range("D45").Value2 = "'41.78"
$ws.Range("E45").Value2 = "  +0.05%  "
$ws.Range("B46").Value2 = "Filecoin"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value2 = "'4.37"
$ws.Range("E46").Value2 = "  +0.61%  "
$ws.Range("B47").Value2 = "ONDO"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value2 = "'1.20"
$ws.Range("E47").Value2 = "  -0.21%  "
$ws.Range("E48").Value2 = "  +1.29%  "
$ws.Range("D49").Value2 = "2.465.98"
$ws.Range("E49").Value2 = "  +1.21%  "
$ws.Range("D50").Value2 = "'6.73"
$ws.Range("E50").Value2 = "  -1.95%  "
$ws.Range("D51").Value2 = "'0.898"
$ws.Range("E51").Value2 = "  +0.96%  "

Write-Host "Applied 98 cell updates to cryptos sheet"
